# cambios de las fracciones
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the quarterly period dates (Q3 2022 -> Q4 2022)
$ws.Range("B8").Value = 44835   # 2022-10-01 start of period
$ws.Range("C8").Value = 44926   # 2022-12-31 end of period
$ws.Range("E8").Value = 44926   # 2022-12-31 fecha de elaboracion
$ws.Range("H8").Value = 44936   # 2023-01-10 fecha de validacion
$ws.Range("I8").Value = 44936   # 2023-01-10 fecha de actualizacion

# Update the sheet view (scroll position + selection)
$ws.Range("G15").Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 2
